$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2115286666666667
$ws.Range("H2").Value = 0.634586
$ws.Range("I2").Value = 0.08153347995807345
$ws.Range("J2").Value = 0.08153347995807345
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.451416666666667
$ws.Range("N2").Value = 7.35425
$ws.Range("O2").Value = 0.2191928499183569
$ws.Range("P2").Value = 0.2191928499183569
$ws.Range("Q2").Value = 0.5185448989444444
$ws.Range("R2").Value = 4.6669040905
$ws.Range("S2").Value = 0.01787155583577135
$ws.Range("T2").Value = 0.01787155583577135
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2115286666666667
$ws.Range("H3").Value = 0.634586
$ws.Range("I3").Value = 0.08153347995807345
$ws.Range("J3").Value = 0.08153347995807345
$ws.Range("O3").Value = 0.4446889938320204
$ws.Range("P3").Value = 0.4446889938320204
$ws.Range("Q3").Value = 1.052001511245556
$ws.Range("R3").Value = 9.46801360121
$ws.Range("S3").Value = 0.03625704116617888
$ws.Range("T3").Value = 0.03625704116617888
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2115286666666667
$ws.Range("H4").Value = 0.634586
$ws.Range("I4").Value = 0.08153347995807345
$ws.Range("J4").Value = 0.08153347995807345
$ws.Range("O4").Value = 0.3361181562496228
$ws.Range("P4").Value = 0.3361181562496228
$ws.Range("Q4").Value = 0.79515529558
$ws.Range("R4").Value = 7.15639766022
$ws.Range("S4").Value = 0.02740488295612322
$ws.Range("T4").Value = 0.02740488295612322
$ws.Range("I5").Value = 0.4735790235655714
$ws.Range("J5").Value = 0.4735790235655714
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.451416666666667
$ws.Range("N5").Value = 7.35425
$ws.Range("O5").Value = 0.2191928499183569
$ws.Range("P5").Value = 0.2191928499183569
$ws.Range("Q5").Value = 3.011915927583333
$ws.Range("R5").Value = 27.10724334825
$ws.Range("S5").Value = 0.1038051358368903
$ws.Range("T5").Value = 0.1038051358368903
$ws.Range("I6").Value = 0.4735790235655714
$ws.Range("J6").Value = 0.4735790235655714
$ws.Range("O6").Value = 0.4446889938320204
$ws.Range("P6").Value = 0.4446889938320204
$ws.Range("S6").Value = 0.2105953794893246
$ws.Range("T6").Value = 0.2105953794893246
$ws.Range("I7").Value = 0.4735790235655714
$ws.Range("J7").Value = 0.4735790235655714
$ws.Range("O7").Value = 0.3361181562496228
$ws.Range("P7").Value = 0.3361181562496228
$ws.Range("S7").Value = 0.1591785082393565
$ws.Range("T7").Value = 0.1591785082393565
$ws.Range("I8").Value = 0.4448874964763552
$ws.Range("J8").Value = 0.4448874964763552
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.451416666666667
$ws.Range("N8").Value = 7.35425
$ws.Range("O8").Value = 0.2191928499183569
$ws.Range("P8").Value = 0.2191928499183569
$ws.Range("Q8").Value = 2.829440642305556
$ws.Range("R8").Value = 25.46496578075
$ws.Range("S8").Value = 0.09751615824569526
$ws.Range("T8").Value = 0.09751615824569525
$ws.Range("I9").Value = 0.4448874964763552
$ws.Range("J9").Value = 0.4448874964763552
$ws.Range("O9").Value = 0.4446889938320204
$ws.Range("P9").Value = 0.4446889938320204
$ws.Range("S9").Value = 0.1978365731765169
$ws.Range("T9").Value = 0.1978365731765169
$ws.Range("I10").Value = 0.4448874964763552
$ws.Range("J10").Value = 0.4448874964763552
$ws.Range("O10").Value = 0.3361181562496228
$ws.Range("P10").Value = 0.3361181562496228
$ws.Range("S10").Value = 0.1495347650541431
$ws.Range("T10").Value = 0.149534765054143
